# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The worker/period detail table (B16:G29) is re-sorted: rows are grouped by
# worker (EVER EDUARDO CARDENAS DE LA OSSA first, then FREDY HUERTAS LOPEZ),
# and within each worker the periods now run in descending order
# (2411 -> 2405) instead of the previous ascending/interleaved order.
# The set of (worker, period, valores) rows is identical - only the order
# (and therefore which row each combination sits in) changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("CC", "1053004704", "EVER EDUARDO CARDENAS DE LA OSSA", "2411", 27402,  1468000),
    @("CC", "1053004704", "EVER EDUARDO CARDENAS DE LA OSSA", "2410", 58720,  1468000),
    @("CC", "1053004704", "EVER EDUARDO CARDENAS DE LA OSSA", "2409", 58720,  1468000),
    @("CC", "1053004704", "EVER EDUARDO CARDENAS DE LA OSSA", "2408", 58720,  1468000),
    @("CC", "1053004704", "EVER EDUARDO CARDENAS DE LA OSSA", "2407", 58720,  1468000),
    @("CC", "1053004704", "EVER EDUARDO CARDENAS DE LA OSSA", "2406", 58720,  1468000),
    @("CC", "1053004704", "EVER EDUARDO CARDENAS DE LA OSSA", "2405", 58720,  1468000),
    @("CC", "74339257",   "FREDY HUERTAS LOPEZ",              "2411", 44000,  2357150),
    @("CC", "74339257",   "FREDY HUERTAS LOPEZ",              "2410", 94286,  2357150),
    @("CC", "74339257",   "FREDY HUERTAS LOPEZ",              "2409", 94286,  2357150),
    @("CC", "74339257",   "FREDY HUERTAS LOPEZ",              "2408", 94286,  2357150),
    @("CC", "74339257",   "FREDY HUERTAS LOPEZ",              "2407", 94286,  2357150),
    @("CC", "74339257",   "FREDY HUERTAS LOPEZ",              "2406", 94286,  2357150),
    @("CC", "74339257",   "FREDY HUERTAS LOPEZ",              "2405", 94286,  2357150)
)

$startRow = 16
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
    $ws.Cells.Item($row, 6).Value = $rec[4]
    $ws.Cells.Item($row, 7).Value = $rec[5]
}
